$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / header text updates (shared strings) ---
$ws.Range("A8").Value = "Volume 31   Number  35"
$ws.Range("C9").Value = "Report Covering the Week  8/26/2024  Through  9/1/2024"

# --- Simple numeric value updates (style unchanged) ---
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = -84
$ws.Range("G15").Value = 2
$ws.Range("J15").Value = 12
$ws.Range("K15").Value = 16.666666666666
$ws.Range("N15").Value = -75
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -75
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = -46.153846153846
$ws.Range("I16").Value = 55
$ws.Range("J16").Value = 89
$ws.Range("K16").Value = -38.202247191011
$ws.Range("L16").Value = -52.991452991453
$ws.Range("M16").Value = -68.926553672316
$ws.Range("N16").Value = -95.119787045252
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = -61.538461538461
$ws.Range("F17").Value = 23
$ws.Range("H17").Value = -25.806451612903
$ws.Range("I17").Value = 174
$ws.Range("J17").Value = 231
$ws.Range("K17").Value = -24.675324675324
$ws.Range("L17").Value = -24.347826086956
$ws.Range("M17").Value = -14.285714285714
$ws.Range("N17").Value = -72.727272727272
$ws.Range("D18").Value = 7
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = -69.230769230769
$ws.Range("I18").Value = 46
$ws.Range("J18").Value = 55
$ws.Range("K18").Value = -16.363636363636
$ws.Range("L18").Value = -48.888888888888
$ws.Range("M18").Value = -77.560975609756
$ws.Range("N18").Value = -96.372239747634
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = -40
$ws.Range("F19").Value = 29
$ws.Range("G19").Value = 39
$ws.Range("H19").Value = -25.641025641025
$ws.Range("I19").Value = 165
$ws.Range("J19").Value = 287
$ws.Range("K19").Value = -42.508710801393
$ws.Range("L19").Value = -30.962343096234
$ws.Range("M19").Value = -38.888888888888
$ws.Range("N19").Value = -63.251670378619
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("I20").Value = 52
$ws.Range("J20").Value = 69
$ws.Range("K20").Value = -24.637681159420
$ws.Range("L20").Value = -45.833333333333
$ws.Range("M20").Value = -50.943396226415
$ws.Range("N20").Value = -94.196428571428
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 41
$ws.Range("E21").Value = -60.975609756097
$ws.Range("F21").Value = 69
$ws.Range("G21").Value = 105
$ws.Range("H21").Value = -34.285714285714
$ws.Range("I21").Value = 510
$ws.Range("J21").Value = 746
$ws.Range("K21").Value = -31.635388739946
$ws.Range("L21").Value = -35.031847133758
$ws.Range("M21").Value = -47.638603696098
$ws.Range("N21").Value = -88.562457950213
$ws.Range("I22").Value = 8
$ws.Range("K22").Value = -33.333333333333
$ws.Range("L22").Value = -46.666666666666
$ws.Range("M22").Value = -75.757575757575
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = -16.666666666666
$ws.Range("F24").Value = 73
$ws.Range("G24").Value = 104
$ws.Range("H24").Value = -29.807692307692
$ws.Range("I24").Value = 583
$ws.Range("J24").Value = 754
$ws.Range("K24").Value = -22.679045092838
$ws.Range("L24").Value = -22.37017310253
$ws.Range("M24").Value = -14.011799410029
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = -50
$ws.Range("F25").Value = 14
$ws.Range("G25").Value = 30
$ws.Range("H25").Value = -53.333333333333
$ws.Range("I25").Value = 105
$ws.Range("J25").Value = 165
$ws.Range("K25").Value = -36.363636363636
$ws.Range("L25").Value = -7.079646017699
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = -20
$ws.Range("F26").Value = 33
$ws.Range("G26").Value = 42
$ws.Range("H26").Value = -21.428571428571
$ws.Range("I26").Value = 369
$ws.Range("J26").Value = 340
$ws.Range("K26").Value = 8.529411764705
$ws.Range("L26").Value = -7.518796992481
$ws.Range("M26").Value = -32.786885245901
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -66.666666666666
$ws.Range("J27").Value = 15
$ws.Range("K27").Value = 20
$ws.Range("C28").Value = 4
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 166.666666666667
$ws.Range("I28").Value = 54
$ws.Range("K28").Value = 42.105263157894
$ws.Range("L28").Value = 35
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 100
$ws.Range("I29").Value = 12
$ws.Range("K29").Value = 20
$ws.Range("L29").Value = -29.411764705882
$ws.Range("M29").Value = -55.555555555555
$ws.Range("N29").Value = -85.185185185185
$ws.Range("F30").Value = 2
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 100
$ws.Range("I30").Value = 10
$ws.Range("K30").Value = 11.111111111111
$ws.Range("L30").Value = -33.333333333333
$ws.Range("M30").Value = -52.380952380952
$ws.Range("N30").Value = -87.5
$ws.Range("F31").Value = 6
$ws.Range("H31").Value = 500

# --- Updates that also require a style/type change (text<->number) ---
# Set the value FIRST (using a leading apostrophe for numeric-looking text so Excel
# does not auto-coerce it back to a number), THEN copy/paste-special the formats from
# a stable donor cell of the target style so the final style matches exactly and the
# quote-prefix bookkeeping from the text entry does not linger as a new style.
$ws.Range("D15").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D15").PasteSpecial(-4122)

$ws.Range("E15").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E15").PasteSpecial(-4122)

$ws.Range("C16").Value = 1
$ws.Range("I14").Copy()
$ws.Range("C16").PasteSpecial(-4122)

$ws.Range("C22").Value = 1
$ws.Range("I14").Copy()
$ws.Range("C22").PasteSpecial(-4122)

$ws.Range("F22").Value = 1
$ws.Range("I14").Copy()
$ws.Range("F22").PasteSpecial(-4122)

$ws.Range("D27").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D27").PasteSpecial(-4122)

$ws.Range("E27").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E27").PasteSpecial(-4122)

$ws.Range("D28").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("D28").PasteSpecial(-4122)

$ws.Range("E28").Value = "'***.*"
$ws.Range("A14").Copy()
$ws.Range("E28").PasteSpecial(-4122)

$ws.Range("C29").Value = 1
$ws.Range("I14").Copy()
$ws.Range("C29").PasteSpecial(-4122)

$ws.Range("C30").Value = 1
$ws.Range("I14").Copy()
$ws.Range("C30").PasteSpecial(-4122)

$excel.CutCopyMode = 0
Write-Host "Edit complete"
